$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New journal entries (rows 20-23). Within each row the historical entry
# order was: Activité, Date, Description, temps (the time spent was noted
# after writing the description), which is what produced the shared-string
# ordering in the saved file, so we replay writes in that same order.

$rows = @(
    @{ Row = 20; Activite = "Formation Laravel"; Date = 43509; Temps = "3h00"; Description = "J'ai commencé à approché les bases de données avec laravel" },
    @{ Row = 21; Activite = "Analyse du code de RackaKey"; Date = 43510; Temps = "2h30"; Description = "Analyse de code du projet RackaKey" },
    @{ Row = 22; Activite = "Laravel et BDD"; Date = 43515; Temps = "5h00"; Description = "J'ai commencé à apprendre la façon dont laravel gère les bases de données, c'est assez complexe à comprendre" },
    @{ Row = 23; Activite = "Laravel et BDD"; Date = 43516; Temps = "1h00"; Description = "Je continue à apprendre comment utiliser les bases de données avec laravel" }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    $ws.Cells.Item($rowIndex, 1).Value = $r.Activite

    $ws.Cells.Item($rowIndex, 2).Value = $r.Date
    # Copy the date number format from the row above so the new date cell
    # reuses the existing date style instead of creating a new one.
    $ws.Cells.Item($rowIndex - 1, 2).Copy()
    $ws.Cells.Item($rowIndex, 2).PasteSpecial(-4122)
    $ws.Cells.Item($rowIndex, 2).Value = $r.Date

    $ws.Cells.Item($rowIndex, 4).Value = $r.Description
    $ws.Cells.Item($rowIndex, 3).Value = $r.Temps
}

$excel.CutCopyMode = 0

# Column widths were tweaked as part of this update (target stored widths
# are 27 / 12.140625 / 7.85546875 characters; the COM ColumnWidth setter
# quantizes to a 1/6-character grid, so we pick the input that lands on the
# nearest achievable grid point to each target).
$ws.Columns.Item(1).ColumnWidth = 26.166666666666668
$ws.Columns.Item(2).ColumnWidth = 11.333333333333334
$ws.Columns.Item(3).ColumnWidth = 7.0

# Move the active selection to the next empty row, like Excel does after
# the last row of data has been entered.
$ws.Range("A24").Select()
